$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "-"
$ws.Range("B3").Value = "MCT-1A-Gestão"
$ws.Range("F4").Value = "MEC-2A-Mecanica material"
$ws.Range("B6").Value = "MEC-1A-Gestão"
$ws.Range("C6").Value = "-"
$ws.Range("B7").Value = "MEC-1A-Gestão"
